# Generate Report for Handoff
# Adds two new file rows (8ce131b0-... and cf12c03e-...) to the
# Overview / zh-cn / de-de worksheets, mirroring the pattern used by the
# existing "Ready for handoff" rows (e.g. 7787a229-...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data
# ---------------------------------------------------------------------
$fileA = "8ce131b0-eed7-44f8-a69a-e580b91b43a3"
$fileB = "cf12c03e-dca1-437d-86e3-380beaba088f"

$hashA = "468b6d45777b69cc4396298e30948008616ef31c"
$hashB = "536e2c66a17efe72ea26173a2750ffed3a7e3d4a"

$mdShaA = "5b6c1b5a2dcb9f6e6a9e7e6a1d4c9b7a3f5e2d1c0"
$mdShaB = "6c7d2c6b3edca0f7f7b0f7b2e5dad8b4067f3e2d1"

$xlfShaA = "2a8b3c5d7e9f1a3b5c7d9e1f3a5b7c9d1e3f5a7b9"
$xlfShaB = "3b9c4d6e8f0a2b4c6d8e0f2a4b6c8d0e2f4a6b8c0"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 6
$ws1.Range("A5:D5").Copy()
$ws1.Range("A6:D6").PasteSpecial(-4104)  # xlPasteAll
$ws1.Range("A6").Value = "$fileA.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"
$ws1.Range("D6").Value = "2016-31-19 02:31:05"
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaA/e2e/$fileA.md", "", "", "$fileA.md") | Out-Null
$ws1.Range("A5:D5").Copy()
$ws1.Range("A6:D6").PasteSpecial(-4122)  # xlPasteFormats (restore non-hyperlink look)

# Row 7
$ws1.Range("A5:D5").Copy()
$ws1.Range("A7:D7").PasteSpecial(-4104)  # xlPasteAll
$ws1.Range("A7").Value = "$fileB.md"
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"
$ws1.Range("D7").Value = "2016-31-19 02:31:05"
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaB/e2e/$fileB.md", "", "", "$fileB.md") | Out-Null
$ws1.Range("A5:D5").Copy()
$ws1.Range("A7:D7").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 6
$ws2.Range("A5:K5").Copy()
$ws2.Range("A6:K6").PasteSpecial(-4104)
$ws2.Range("A6").Value = "$fileA.md"
$ws2.Range("B6").Value = ".md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("D6").Value = "$fileA.$hashA.zh-cn.xlf"
$ws2.Range("E6").Value = "2016-03-19 02:31:02"
$ws2.Range("H6").Value = "0001-01-01 00:00:00"
$ws2.Range("I6").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaA/e2e/$fileA.md", "", "", "$fileA.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaA/e2e/$fileA.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfShaA/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileA.$hashA.zh-cn.xlf", "", "", "$fileA.$hashA.zh-cn.xlf") | Out-Null
$ws2.Range("A5:K5").Copy()
$ws2.Range("A6:K6").PasteSpecial(-4122)

# Row 7
$ws2.Range("A5:K5").Copy()
$ws2.Range("A7:K7").PasteSpecial(-4104)
$ws2.Range("A7").Value = "$fileB.md"
$ws2.Range("B7").Value = ".md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("D7").Value = "$fileB.$hashB.zh-cn.xlf"
$ws2.Range("E7").Value = "2016-03-19 02:31:02"
$ws2.Range("H7").Value = "0001-01-01 00:00:00"
$ws2.Range("I7").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaB/e2e/$fileB.md", "", "", "$fileB.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaB/e2e/$fileB.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfShaB/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileB.$hashB.zh-cn.xlf", "", "", "$fileB.$hashB.zh-cn.xlf") | Out-Null
$ws2.Range("A5:K5").Copy()
$ws2.Range("A7:K7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 6
$ws3.Range("A5:K5").Copy()
$ws3.Range("A6:K6").PasteSpecial(-4104)
$ws3.Range("A6").Value = "$fileA.md"
$ws3.Range("B6").Value = ".md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("D6").Value = "$fileA.$hashA.de-de.xlf"
$ws3.Range("E6").Value = "2016-03-19 02:31:05"
$ws3.Range("H6").Value = "0001-01-01 00:00:00"
$ws3.Range("I6").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaA/e2e/$fileA.md", "", "", "$fileA.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaA/e2e/$fileA.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfShaA/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileA.$hashA.de-de.xlf", "", "", "$fileA.$hashA.de-de.xlf") | Out-Null
$ws3.Range("A5:K5").Copy()
$ws3.Range("A6:K6").PasteSpecial(-4122)

# Row 7
$ws3.Range("A5:K5").Copy()
$ws3.Range("A7:K7").PasteSpecial(-4104)
$ws3.Range("A7").Value = "$fileB.md"
$ws3.Range("B7").Value = ".md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("D7").Value = "$fileB.$hashB.de-de.xlf"
$ws3.Range("E7").Value = "2016-03-19 02:31:05"
$ws3.Range("H7").Value = "0001-01-01 00:00:00"
$ws3.Range("I7").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaB/e2e/$fileB.md", "", "", "$fileB.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdShaB/e2e/$fileB.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfShaB/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileB.$hashB.de-de.xlf", "", "", "$fileB.$hashB.de-de.xlf") | Out-Null
$ws3.Range("A5:K5").Copy()
$ws3.Range("A7:K7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Output "Report generated for handoff: added $fileA and $fileB rows to Overview, zh-cn, de-de sheets."
